# Sankey: AMPLpy + color renovation
# Header row was reordered (target now first, source second) and the
# new "target" header is highlighted with an explicit black font.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap A1 and B1 header values: A1 becomes "target", B1 becomes "source"
$ws.Range("A1").Value = "target"
$ws.Range("B1").Value = "source"

# Give the new A1 header an explicit black font color
$ws.Range("A1").Font.Color = 0

# Move the active selection to A2, matching the saved workbook view
$ws.Range("A2").Select()
